$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 with the new value
$ws.Range("A2").Value = "M6000_V73"

# Remove the now-obsolete rows 3 through 6 so the sheet dimension shrinks to A1:A2
$ws.Rows("3:6").Delete()
